$wb = $excel.ActiveWorkbook

# Rename the existing sheet "t0" -> "t1_pre"
$ws1 = $wb.ActiveSheet
$ws1.Name = "t1_pre"

# Add a new sheet "t1_post" right after "t1_pre"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "t1_post"

# --- Populate t1_post with the balance-check table ---
$ws2.Range("A1").Value = 'var'
$ws2.Range("B1").Value = 'level'
$ws2.Range("C1").Value = 'ctrl.mean/n'
$ws2.Range("D1").Value = 'ctrl.sd/%'
$ws2.Range("E1").Value = 'ONE.mean/n'
$ws2.Range("F1").Value = 'ONE.sd/%'
$ws2.Range("G1").Value = 'p.val'

$ws2.Range("A2").Value = 'n'
$ws2.Range("C2").Value = "'302"
$ws2.Range("D2").Value = 100
$ws2.Range("E2").Value = "'151"
$ws2.Range("F2").Value = 100

$ws2.Range("A3").Value = 'index'
$ws2.Range("C3").Value = "'2017-05-08"
$ws2.Range("D3").Value = 203.07
$ws2.Range("E3").Value = "'2018-10-15"
$ws2.Range("F3").Value = 41.54
$ws2.Range("G3").Value = "'0.000"

$ws2.Range("A4").Value = 'male'
$ws2.Range("B4").Value = 1
$ws2.Range("C4").Value = "'94"
$ws2.Range("D4").Value = 31.1
$ws2.Range("E4").Value = "'48"
$ws2.Range("F4").Value = 31.8
$ws2.Range("G4").Value = "'0.971"

$ws2.Range("A5").Value = 'age'
$ws2.Range("C5").Value = "'43.28"
$ws2.Range("D5").Value = 10.57
$ws2.Range("E5").Value = "'42.72"
$ws2.Range("F5").Value = 11.11
$ws2.Range("G5").Value = "'0.610"

$ws2.Range("A6").Value = 'unempl'
$ws2.Range("C6").Value = "'4.23"
$ws2.Range("D6").Value = 0.11
$ws2.Range("E6").Value = "'3.82"
$ws2.Range("F6").Value = 0.04
$ws2.Range("G6").Value = "'0.000"

$ws2.Range("A7").Value = 'b_46'
$ws2.Range("B7").Value = 1
$ws2.Range("C7").Value = "'20"
$ws2.Range("D7").Value = 6.6
$ws2.Range("E7").Value = "'13"
$ws2.Range("F7").Value = 8.6
$ws2.Range("G7").Value = "'0.565"

$ws2.Range("A8").Value = 'b_47'
$ws2.Range("B8").Value = 1
$ws2.Range("C8").Value = "'30"
$ws2.Range("D8").Value = 9.9
$ws2.Range("E8").Value = "'20"
$ws2.Range("F8").Value = 13.2
$ws2.Range("G8").Value = "'0.367"

$ws2.Range("A9").Value = 'b_84'
$ws2.Range("B9").Value = 1
$ws2.Range("C9").Value = "'20"
$ws2.Range("D9").Value = 6.6
$ws2.Range("E9").Value = "'12"
$ws2.Range("F9").Value = 7.9
$ws2.Range("G9").Value = "'0.746"

$ws2.Range("A10").Value = 'b_85'
$ws2.Range("B10").Value = 1
$ws2.Range("C10").Value = "'43"
$ws2.Range("D10").Value = 14.2
$ws2.Range("E10").Value = "'21"
$ws2.Range("F10").Value = 13.9
$ws2.Range("G10").Value = "'1.000"

$ws2.Range("A11").Value = 'b_86'
$ws2.Range("B11").Value = 1
$ws2.Range("C11").Value = "'17"
$ws2.Range("D11").Value = 5.6
$ws2.Range("E11").Value = "'6"
$ws2.Range("F11").Value = 4
$ws2.Range("G11").Value = "'0.596"

$ws2.Range("A12").Value = 'b_87'
$ws2.Range("B12").Value = 1
$ws2.Range("C12").Value = "'26"
$ws2.Range("D12").Value = 8.6
$ws2.Range("E12").Value = "'13"
$ws2.Range("F12").Value = 8.6
$ws2.Range("G12").Value = "'1.000"

$ws2.Range("A13").Value = 'b_88'
$ws2.Range("B13").Value = 1
$ws2.Range("C13").Value = "'45"
$ws2.Range("D13").Value = 14.9
$ws2.Range("E13").Value = "'22"
$ws2.Range("F13").Value = 14.6
$ws2.Range("G13").Value = "'1.000"

$ws2.Range("A14").Value = 'b_99'
$ws2.Range("B14").Value = 1
$ws2.Range("C14").Value = "'20"
$ws2.Range("D14").Value = 6.6
$ws2.Range("E14").Value = "'9"
$ws2.Range("F14").Value = 6
$ws2.Range("G14").Value = "'0.946"

$ws2.Range("A15").Value = 'y_none'
$ws2.Range("C15").Value = "'73.97"
$ws2.Range("D15").Value = 26.29
$ws2.Range("E15").Value = "'71.69"
$ws2.Range("F15").Value = 28.07
$ws2.Range("G15").Value = "'0.392"

$ws2.Range("A16").Value = 'y_dgp'
$ws2.Range("C16").Value = "'7.92"
$ws2.Range("D16").Value = 17.44
$ws2.Range("E16").Value = "'8.99"
$ws2.Range("F16").Value = 19.41
$ws2.Range("G16").Value = "'0.974"

$ws2.Range("A17").Value = 'y_edu'
$ws2.Range("C17").Value = "'0.82"
$ws2.Range("D17").Value = 5.62
$ws2.Range("E17").Value = "'0.83"
$ws2.Range("F17").Value = 4.68
$ws2.Range("G17").Value = "'0.723"

$ws2.Range("A18").Value = 'y_sgdp'
$ws2.Range("C18").Value = "'15.39"
$ws2.Range("D18").Value = 11.17
$ws2.Range("E18").Value = "'16.54"
$ws2.Range("F18").Value = 10.55
$ws2.Range("G18").Value = "'0.022"

$ws2.Range("A19").Value = 'y_baby'
$ws2.Range("C19").Value = "'3.19"
$ws2.Range("D19").Value = 11.86
$ws2.Range("E19").Value = "'3.26"
$ws2.Range("F19").Value = 12.2
$ws2.Range("G19").Value = "'0.689"

$ws2.Range("A20").Value = 'y_flex'
$ws2.Range("C20").Value = "'0.53"
$ws2.Range("D20").Value = 6.05
$ws2.Range("E20").Value = "'0.53"
$ws2.Range("F20").Value = 6.51
$ws2.Range("G20").Value = "'0.730"

$ws2.Range("A21").Value = 'y_cash'
$ws2.Range("C21").Value = "'0.34"
$ws2.Range("D21").Value = 3.22
$ws2.Range("E21").Value = "'0.50"
$ws2.Range("F21").Value = 3.67
$ws2.Range("G21").Value = "'0.650"

$ws2.Range("A22").Value = 'y_reval'
$ws2.Range("C22").Value = "'0.00"
$ws2.Range("D22").Value = 0
$ws2.Range("E22").Value = "'0.00"
$ws2.Range("F22").Value = 0
$ws2.Range("G22").Value = 'NaN'

# Header row formatting to match the style used in t1_pre
$header = $ws2.Range("A1:G1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108  # xlCenter

$ws1.Activate()
